$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "JD_002"
$ws.Range("B3").Value = "Senior Fullstack Engineer"
$ws.Range("C3").Value = "We are seeking a Junior RPA Developer to design, develop, and support automation solutions.
Collaborate with teams to streamline business processes using RPA tools like UiPath or Automation Anywhere. Join Akkodis to grow your skills in a dynamic, tech-driven environment"
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0

$ws.Rows(3).AutoFit()
